# Scheduled-runner price/profit refresh across the per-job sheets (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR). Each block below updates the
# currentAveragePrice* / LevePrice* / LeveProfit* columns (H..N) for the
# specific rows whose market data moved since the last run. Where a row had
# no HQ-profit (M) value before, it is newly populated; where a row's NQ
# craft became fully unprofitable/unused, N is cleared instead of written.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 30624
$ws.Range("J95").Value = 30624
$ws.Range("L95").Value = 30624
$ws.Range("N95").Value = -36116
$ws.Range("H100").Value = 1885.1818
$ws.Range("I100").Value = 1376.7142
$ws.Range("J100").Value = 2775
$ws.Range("K100").Value = 1376.7142
$ws.Range("L100").Value = 2775
$ws.Range("M100").Value = -835.7141999999999
$ws.Range("N100").Value = -3857
$ws.Range("H116").Value = 3331.4167
$ws.Range("I116").Value = 2499
$ws.Range("J116").Value = 3407.0908
$ws.Range("K116").Value = 2499
$ws.Range("L116").Value = 3407.0908
$ws.Range("M116").Value = 943
$ws.Range("N116").Value = -10291.0908
$ws.Range("H137").Value = 1554.2195
$ws.Range("I137").Value = 1430.4445
$ws.Range("J137").Value = 1792.9286
$ws.Range("K137").Value = 4291.333500000001
$ws.Range("L137").Value = 5378.7858
$ws.Range("M137").Value = -1741.333500000001
$ws.Range("N137").Value = -10478.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3117.0466
$ws.Range("I32").Value = 1778.3055
$ws.Range("K32").Value = 1778.3055
$ws.Range("M32").Value = -1491.3055
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H61").Value = 8764.444
$ws.Range("I61").Value = 11376
$ws.Range("J61").Value = 5500
$ws.Range("K61").Value = 11376
$ws.Range("L61").Value = 5500
$ws.Range("M61").Value = -11164
$ws.Range("N61").Value = -5924
$ws.Range("H132").Value = 19202.45
$ws.Range("I132").Value = 1409.5
$ws.Range("K132").Value = 4228.5
$ws.Range("M132").Value = -1698.5
$ws.Range("H133").Value = 39999
$ws.Range("J133").Value = 39999
$ws.Range("L133").Value = 39999
$ws.Range("N133").Value = -45059
$ws.Range("H136").Value = 8764.444
$ws.Range("I136").Value = 11376
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 34128
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -31578
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").Value = ""
$ws.Range("H95").Value = 26324.6
$ws.Range("J95").Value = 26324.6
$ws.Range("L95").Value = 26324.6
$ws.Range("N95").Value = -31816.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2264.838
$ws.Range("I31").Value = 1002.76
$ws.Range("J31").Value = 4894.1665
$ws.Range("K31").Value = 1002.76
$ws.Range("L31").Value = 4894.1665
$ws.Range("M31").Value = -707.76
$ws.Range("N31").Value = -5484.1665
$ws.Range("H34").Value = 2264.838
$ws.Range("I34").Value = 1002.76
$ws.Range("J34").Value = 4894.1665
$ws.Range("K34").Value = 1002.76
$ws.Range("L34").Value = 4894.1665
$ws.Range("M34").Value = -800.76
$ws.Range("N34").Value = -5298.1665
$ws.Range("H105").Value = 11364799
$ws.Range("I105").Value = 13889665
$ws.Range("J105").Value = 2900
$ws.Range("K105").Value = 13889665
$ws.Range("L105").Value = 2900
$ws.Range("M105").Value = -13887918
$ws.Range("N105").Value = -6394
$ws.Range("H132").Value = 3695.5557
$ws.Range("I132").Value = 1173.75
$ws.Range("J132").Value = 5713
$ws.Range("K132").Value = 3521.25
$ws.Range("L132").Value = 17139
$ws.Range("M132").Value = -991.25
$ws.Range("N132").Value = -22199

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 505.5
$ws.Range("J23").Value = 684.1818
$ws.Range("L23").Value = 2052.5454
$ws.Range("N23").Value = -2522.5454
$ws.Range("H68").Value = 1139.5
$ws.Range("J68").Value = 1143.7273
$ws.Range("L68").Value = 3431.1819
$ws.Range("N68").Value = -5053.1819
$ws.Range("H71").Value = 1139.5
$ws.Range("J71").Value = 1143.7273
$ws.Range("L71").Value = 10293.5457
$ws.Range("N71").Value = -18405.5457
$ws.Range("H122").Value = 1127
$ws.Range("I122").Value = 504
$ws.Range("J122").Value = 1334.6666
$ws.Range("K122").Value = 4536
$ws.Range("L122").Value = 12011.9994
$ws.Range("M122").Value = -2086
$ws.Range("N122").Value = -16911.9994
$ws.Range("H131").Value = 800.64
$ws.Range("I131").Value = 388
$ws.Range("J131").Value = 822.3579
$ws.Range("K131").Value = 1164
$ws.Range("L131").Value = 2467.0737
$ws.Range("M131").Value = 3876
$ws.Range("N131").Value = -12547.0737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 27750
$ws.Range("J39").Value = 27750
$ws.Range("L39").Value = 27750
$ws.Range("N39").Value = -28814
$ws.Range("H80").Value = 3580.25
$ws.Range("J80").Value = 4160
$ws.Range("L80").Value = 4160
$ws.Range("N80").Value = -6156
$ws.Range("H83").Value = 3580.25
$ws.Range("J83").Value = 4160
$ws.Range("L83").Value = 20800
$ws.Range("N83").Value = -30784
$ws.Range("H97").Value = 1500.8928
$ws.Range("J97").Value = 3022.111
$ws.Range("L97").Value = 3022.111
$ws.Range("N97").Value = -4014.111
$ws.Range("H107").Value = 1856.5714
$ws.Range("I107").Value = 419.2
$ws.Range("J107").Value = 5450
$ws.Range("K107").Value = 419.2
$ws.Range("L107").Value = 5450
$ws.Range("M107").Value = 1500.8
$ws.Range("N107").Value = -9290
$ws.Range("H132").Value = 43479.54
$ws.Range("I132").Value = 5248.4443
$ws.Range("K132").Value = 15745.3329
$ws.Range("M132").Value = -13215.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 226.54167
$ws.Range("I55").Value = 155.54546
$ws.Range("J55").Value = 286.6154
$ws.Range("K55").Value = 155.54546
$ws.Range("L55").Value = 286.6154
$ws.Range("M55").Value = 17.45454000000001
$ws.Range("N55").Value = -632.6154
$ws.Range("H97").Value = 22497.5
$ws.Range("J97").Value = 22497.5
$ws.Range("L97").Value = 22497.5
$ws.Range("N97").Value = -24479.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 33372
$ws.Range("J95").Value = 33372
$ws.Range("L95").Value = 33372
$ws.Range("N95").Value = -38864
$ws.Range("H100").Value = 1143.2142
$ws.Range("I100").Value = 808.6667
$ws.Range("J100").Value = 1394.125
$ws.Range("K100").Value = 1617.3334
$ws.Range("L100").Value = 2788.25
$ws.Range("M100").Value = -1076.3334
$ws.Range("N100").Value = -3870.25
